$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new "ECs" sending-cluster row (was not present before; replaces the
# old row 2 position in-place since the overall row count drops 5 -> 4)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1214023333333333
$ws.Range("H2").Value = 0.364207
$ws.Range("I2").Value = 0.2856182748266287
$ws.Range("J2").Value = 0.2856182748266287
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9721403333333333
$ws.Range("N2").Value = 2.916421
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1180201047941111
$ws.Range("R2").Value = 1.062180943147
$ws.Range("S2").Value = 0.2856182748266287
$ws.Range("T2").Value = 0.2856182748266287

# Row 3: FAPs -> FAPs (re-run with new TPM values; target cluster flips from
# MuSCs to FAPs and several derived-specificity figures are recomputed)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.174539
$ws.Range("H3").Value = 0.523617
$ws.Range("I3").Value = 0.4106307243130825
$ws.Range("J3").Value = 0.4106307243130825
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9721403333333333
$ws.Range("N3").Value = 2.916421
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.1696764016396667
$ws.Range("R3").Value = 1.527087614757
$ws.Range("S3").Value = 0.4106307243130825
$ws.Range("T3").Value = 0.4106307243130825

# Row 4: MuSCs -> FAPs (re-run with new TPM values)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1291096666666667
$ws.Range("H4").Value = 0.387329
$ws.Range("I4").Value = 0.3037510008602889
$ws.Range("J4").Value = 0.3037510008602889
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9721403333333333
$ws.Range("N4").Value = 2.916421
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.1255127143898889
$ws.Range("R4").Value = 1.129614429509
$ws.Range("S4").Value = 0.3037510008602889
$ws.Range("T4").Value = 0.3037510008602889

# Row 5 (old MuSCs -> MuSCs pair) no longer exists in the refreshed output
$ws.Rows.Item(5).Delete()
